$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.246.81"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.035.03"
$ws.Range("E3").Value = "  +3.97%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.19"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.69"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  +7.40%  "
$ws.Range("D10").Value = "3.033.83"
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  +6.79%  "
$ws.Range("D14").Value = "3.593.52"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.84"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "76.097.50"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "3.024.05"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.92"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.01"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  +6.06%  "
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.43"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.92"
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "490.46"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +12.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.56"
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.85"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.05"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "192.02"
$ws.Range("E40").Value = "  +7.17%  "
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.782"
$ws.Range("E45").Value = "  +19.86%  "
$ws.Range("E46").Value = "  +6.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.03"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.593"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  -0.27%  "
